$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add two new checklist rows for "Robot Vision" and "Robot LEDs",
# following the same layout/format as the existing rows (e.g. row 37).
$ws.Range("B37:F37").Copy() | Out-Null
$ws.Range("B38:F38").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteAll) | Out-Null
$ws.Range("B39:F39").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteAll) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B38").Value = "Robot Vision"
$ws.Range("C38").Value = "X"
$ws.Range("D38").Formula = '=IF(EXACT(LOWER(C38), "check"), "ü", "û")'
$ws.Range("E38").Value = ""
$ws.Range("F38").Value = ""

# Fix capitalization of the two auto mode names (rows 14 and 15)
$ws.Range("B14").Value = "AutoModeOne_StackTotes"
$ws.Range("B15").Value = "AutoModeTwo_TakeRecycling"

# Rename "Robot" class to "SlideWinder" (row 7)
$ws.Range("B7").Value = "SlideWinder"

$ws.Range("B39").Value = "Robot LEDs"
$ws.Range("C39").Value = "X"
$ws.Range("D39").Formula = '=IF(EXACT(LOWER(C39), "check"), "ü", "û")'
$ws.Range("E39").Value = ""
$ws.Range("F39").Value = ""

# Update the selection to match the author's final cursor position
$ws.Range("E39:F39").Select()
